$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rng = $ws.Range("A2:K2")

# Force text interpretation so date/numeric-looking strings (e.g. "2017-11-05",
# "123456789", "18") are stored as shared strings, not auto-converted to
# Excel dates/numbers, matching the source data row being appended.
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "2017-11-05"
$ws.Range("B2").Value = "18:44:03.520080"
$ws.Range("C2").Value = "dd"
$ws.Range("D2").Value = "123456789"
$ws.Range("E2").Value = "18"
$ws.Range("F2").Value = "m"
$ws.Range("G2").Value = "y"
$ws.Range("H2").Value = "na"
$ws.Range("I2").Value = "y"
$ws.Range("J2").Value = "na"
$ws.Range("K2").Value = "sdp"

# Reset cell formatting back to the default style so the new row matches
# the plain (unstyled) look of the header row.
$rng.Style = "Normal"
